$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CarNameAndPrice")
$ws.Activate()

# Update B3 value from "firefox" to "chrome"
$ws.Range("B3").Value = "chrome"

# Update the selection to D14
$ws.Range("D14").Select()
